$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 867
$ws.Range("I41").Value = 1000.4
$ws.Range("J41").Value = 200
$ws.Range("K41").Value = 1000.4
$ws.Range("L41").Value = 200
$ws.Range("M41").Value = -560.4
$ws.Range("N41").Value = -1080
$ws.Range("H48").Value = 2996
$ws.Range("H56").Value = 2996
$ws.Range("H96").Value = 3524.4
$ws.Range("I96").Value = 1708
$ws.Range("K96").Value = 5124
$ws.Range("M96").Value = -3751
$ws.Range("H98").Value = 1293.8649
$ws.Range("I98").Value = 1255.069
$ws.Range("K98").Value = 1255.069
$ws.Range("M98").Value = 242.931
$ws.Range("H106").Value = 6066.3335
$ws.Range("I106").Value = 3374.5334
$ws.Range("K106").Value = 3374.5334
$ws.Range("M106").Value = -2743.5334
$ws.Range("H107").Value = 1208.25
$ws.Range("I107").Value = 663.5
$ws.Range("J107").Value = 1753
$ws.Range("K107").Value = 663.5
$ws.Range("L107").Value = 1753
$ws.Range("M107").Value = 1256.5
$ws.Range("N107").Value = -5593
$ws.Range("H112").Value = 3198.1333
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3198.1333
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 9594.3999
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -11810.3999
$ws.Range("H122").Value = 1293.8649
$ws.Range("I122").Value = 1255.069
$ws.Range("K122").Value = 3765.207
$ws.Range("M122").Value = -1315.207
$ws.Range("H125").Value = 1767.3
$ws.Range("J125").Value = 1465.5
$ws.Range("L125").Value = 13189.5
$ws.Range("N125").Value = -18109.5
$ws.Range("H132").Value = 6671.0303
$ws.Range("I132").Value = 6946.467
$ws.Range("K132").Value = 20839.401
$ws.Range("M132").Value = -18309.401
$ws.Range("H134").Value = 108570.57
$ws.Range("J134").Value = 108570.57
$ws.Range("L134").Value = 108570.57
$ws.Range("N134").Value = -118710.57
$ws.Range("H138").Value = 3328.4285
$ws.Range("I138").Value = 3366.7407
$ws.Range("K138").Value = 10100.2221
$ws.Range("M138").Value = -4960.222099999999
$ws.Range("H139").Value = 94400
$ws.Range("J139").Value = 94400
$ws.Range("L139").Value = 94400
$ws.Range("N139").Value = -104680
$ws.Range("H140").Value = 126599.664
$ws.Range("J140").Value = 126599.664
$ws.Range("L140").Value = 126599.664
$ws.Range("N140").Value = -136959.664

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4262
$ws.Range("I2").Value = 3946
$ws.Range("K2").Value = 3946
$ws.Range("M2").Value = -3833
$ws.Range("H36").Value = 4511
$ws.Range("J36").Value = 5014.5
$ws.Range("L36").Value = 5014.5
$ws.Range("N36").Value = -5706.5
$ws.Range("H61").Value = 2036.75
$ws.Range("I61").Value = 1952.25
$ws.Range("K61").Value = 1952.25
$ws.Range("M61").Value = -1740.25
$ws.Range("H74").Value = 1485.2903
$ws.Range("I74").Value = 1062.2941
$ws.Range("K74").Value = 1062.2941
$ws.Range("M74").Value = -188.2941000000001
$ws.Range("H77").Value = 1485.2903
$ws.Range("I77").Value = 1062.2941
$ws.Range("K77").Value = 5311.4705
$ws.Range("M77").Value = -943.4705000000004
$ws.Range("H110").Value = 3607.2856
$ws.Range("I110").Value = 3375.1667
$ws.Range("K110").Value = 3375.1667
$ws.Range("M110").Value = -1330.1667
$ws.Range("H116").Value = 4262
$ws.Range("I116").Value = 3946
$ws.Range("K116").Value = 3946
$ws.Range("M116").Value = -1652
$ws.Range("H122").Value = 2002.2727
$ws.Range("I122").Value = 1807.7222
$ws.Range("K122").Value = 5423.1666
$ws.Range("M122").Value = -2973.1666
$ws.Range("H132").Value = 30384.656
$ws.Range("I132").Value = 35185.465
$ws.Range("K132").Value = 105556.395
$ws.Range("M132").Value = -103026.395
$ws.Range("H136").Value = 2036.75
$ws.Range("I136").Value = 1952.25
$ws.Range("K136").Value = 5856.75
$ws.Range("M136").Value = -3306.75

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4262
$ws.Range("I3").Value = 3946
$ws.Range("K3").Value = 3946
$ws.Range("M3").Value = -3832
$ws.Range("H20").Value = 33337668
$ws.Range("I20").Value = 50003504
$ws.Range("J20").Value = 6000
$ws.Range("K20").Value = 50003504
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = -50003257
$ws.Range("N20").Value = -6494
$ws.Range("H58").Value = 60888.5
$ws.Range("I58").Value = 39999
$ws.Range("J58").Value = 81778
$ws.Range("K58").Value = 39999
$ws.Range("L58").Value = 81778
$ws.Range("M58").Value = -39705
$ws.Range("N58").Value = -82366
$ws.Range("H74").Value = 59989
$ws.Range("J74").Value = 59989
$ws.Range("L74").Value = 59989
$ws.Range("N74").Value = -61861
$ws.Range("H77").Value = 59989
$ws.Range("J77").Value = 59989
$ws.Range("L77").Value = 179967
$ws.Range("N77").Value = -189327
$ws.Range("H81").Value = 76999.5
$ws.Range("J81").Value = 76999.5
$ws.Range("L81").Value = 76999.5
$ws.Range("N81").Value = -79121.5
$ws.Range("H84").Value = 76999.5
$ws.Range("J84").Value = 76999.5
$ws.Range("L84").Value = 230998.5
$ws.Range("N84").Value = -241606.5
$ws.Range("H86").Value = 4811.2666
$ws.Range("I86").Value = 3999
$ws.Range("K86").Value = 3999
$ws.Range("M86").Value = -2876
$ws.Range("H89").Value = 4811.2666
$ws.Range("I89").Value = 3999
$ws.Range("K89").Value = 19995
$ws.Range("M89").Value = -14379
$ws.Range("H134").Value = 2179.2258
$ws.Range("I134").Value = 2179.2258
$ws.Range("K134").Value = 6537.6774
$ws.Range("M134").Value = -4002.6774

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1406
$ws.Range("I16").Value = 508.5
$ws.Range("K16").Value = 508.5
$ws.Range("M16").Value = -221.5
$ws.Range("H20").Value = 79696.25
$ws.Range("J20").Value = 79696.25
$ws.Range("L20").Value = 79696.25
$ws.Range("N20").Value = -80168.25
$ws.Range("H30").Value = 79696.25
$ws.Range("J30").Value = 79696.25
$ws.Range("L30").Value = 79696.25
$ws.Range("N30").Value = -79878.25
$ws.Range("H31").Value = 7347.9414
$ws.Range("I31").Value = 7460.2856
$ws.Range("J31").Value = 7269.3
$ws.Range("K31").Value = 7460.2856
$ws.Range("L31").Value = 7269.3
$ws.Range("M31").Value = -7165.2856
$ws.Range("N31").Value = -7859.3
$ws.Range("H34").Value = 7347.9414
$ws.Range("I34").Value = 7460.2856
$ws.Range("J34").Value = 7269.3
$ws.Range("K34").Value = 7460.2856
$ws.Range("L34").Value = 7269.3
$ws.Range("M34").Value = -7258.2856
$ws.Range("N34").Value = -7673.3
$ws.Range("H39").Value = 23512.5
$ws.Range("I39").Value = 18017
$ws.Range("K39").Value = 18017
$ws.Range("M39").Value = -17626
$ws.Range("H49").Value = 23512.5
$ws.Range("I49").Value = 18017
$ws.Range("K49").Value = 18017
$ws.Range("M49").Value = -17835
$ws.Range("H62").Value = 7399
$ws.Range("I62").Value = 7399
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7399
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -6775
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 7399
$ws.Range("I65").Value = 7399
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 36995
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -33875
$ws.Range("N65").ClearContents()
$ws.Range("H86").Value = 9894
$ws.Range("I86").Value = 9894
$ws.Range("K86").Value = 9894
$ws.Range("M86").Value = -8771
$ws.Range("H89").Value = 9894
$ws.Range("I89").Value = 9894
$ws.Range("K89").Value = 49470
$ws.Range("M89").Value = -43854
$ws.Range("H99").Value = 4294.65
$ws.Range("I99").Value = 4318.0835
$ws.Range("J99").Value = 4259.5
$ws.Range("K99").Value = 4318.0835
$ws.Range("L99").Value = 4259.5
$ws.Range("M99").Value = -2820.0835
$ws.Range("N99").Value = -7255.5
$ws.Range("H113").Value = 1406
$ws.Range("I113").Value = 508.5
$ws.Range("K113").Value = 508.5
$ws.Range("M113").Value = 1661.5
$ws.Range("H126").Value = 4294.65
$ws.Range("I126").Value = 4318.0835
$ws.Range("J126").Value = 4259.5
$ws.Range("K126").Value = 12954.2505
$ws.Range("L126").Value = 12778.5
$ws.Range("M126").Value = -10484.2505
$ws.Range("N126").Value = -17718.5
$ws.Range("H128").Value = 79696.25
$ws.Range("J128").Value = 79696.25
$ws.Range("L128").Value = 79696.25
$ws.Range("N128").Value = -89656.25
$ws.Range("H134").Value = 68083.625
$ws.Range("I134").Value = 83333.69500000001
$ws.Range("K134").Value = 250001.085
$ws.Range("M134").Value = -247466.085

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 84.2
$ws.Range("I12").Value = 125.5
$ws.Range("J12").Value = 73.875
$ws.Range("K12").Value = 376.5
$ws.Range("L12").Value = 221.625
$ws.Range("M12").Value = -203.5
$ws.Range("N12").Value = -567.625
$ws.Range("H23").Value = 529
$ws.Range("J23").Value = 529
$ws.Range("L23").Value = 1587
$ws.Range("N23").Value = -2057
$ws.Range("H69").Value = 430
$ws.Range("I69").Value = 200
$ws.Range("K69").Value = 600
$ws.Range("M69").Value = 211
$ws.Range("H72").Value = 430
$ws.Range("I72").Value = 200
$ws.Range("K72").Value = 1800
$ws.Range("M72").Value = 2256

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 390.9524
$ws.Range("I2").Value = 140.72728
$ws.Range("K2").Value = 140.72728
$ws.Range("M2").Value = -27.72728000000001
$ws.Range("H11").Value = 6062.375
$ws.Range("J11").Value = 14499.667
$ws.Range("L11").Value = 14499.667
$ws.Range("N11").Value = -14777.667
$ws.Range("H21").Value = 5000000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 5000000
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H102").Value = 3303.8845
$ws.Range("J102").Value = 4868.6
$ws.Range("L102").Value = 4868.6
$ws.Range("N102").Value = -8112.6
$ws.Range("H122").Value = 2637.5715
$ws.Range("J122").Value = 5295.7144
$ws.Range("L122").Value = 15887.1432
$ws.Range("N122").Value = -20787.1432
$ws.Range("H132").Value = 127765.875
$ws.Range("I132").Value = 145161
$ws.Range("K132").Value = 435483
$ws.Range("M132").Value = -432953

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H61").Value = 7281.4
$ws.Range("I61").Value = 6976.75
$ws.Range("K61").Value = 6976.75
$ws.Range("M61").Value = -6774.75
$ws.Range("H100").Value = 4648.3335
$ws.Range("I100").Value = 4547.857
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 4547.857
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -4006.857
$ws.Range("N100").Value = -6082
$ws.Range("H113").Value = 7281.4
$ws.Range("I113").Value = 6976.75
$ws.Range("K113").Value = 6976.75
$ws.Range("M113").Value = -4806.75
$ws.Range("H122").Value = 4553.1924
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 4941.2104
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 14823.6312
$ws.Range("M122").Value = -8050
$ws.Range("N122").Value = -19723.6312
$ws.Range("H132").Value = 75622.88
$ws.Range("I132").Value = 89699.21000000001
$ws.Range("K132").Value = 269097.63
$ws.Range("M132").Value = -266567.63
$ws.Range("H136").Value = 5778.8
$ws.Range("I136").Value = 3399
$ws.Range("K136").Value = 10197
$ws.Range("M136").Value = -7647

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 22537.5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 22537.5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 22537.5
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -22761.5
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H81").Value = 1472.1875
$ws.Range("J81").Value = 3789
$ws.Range("L81").Value = 7578
$ws.Range("N81").Value = -9700
$ws.Range("H84").Value = 1472.1875
$ws.Range("J84").Value = 3789
$ws.Range("L84").Value = 37890
$ws.Range("N84").Value = -48498
$ws.Range("H107").Value = 1368.6666
$ws.Range("I107").Value = 590.9167
$ws.Range("J107").Value = 2924.1667
$ws.Range("K107").Value = 1772.7501
$ws.Range("L107").Value = 8772.500100000001
$ws.Range("M107").Value = 147.2499
$ws.Range("N107").Value = -12612.5001
$ws.Range("H122").Value = 5019.815
$ws.Range("I122").Value = 5501.5654
$ws.Range("J122").Value = 2249.75
$ws.Range("K122").Value = 16504.6962
$ws.Range("L122").Value = 6749.25
$ws.Range("M122").Value = -14054.6962
$ws.Range("N122").Value = -11649.25
